# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> used by the slide master ("Integral" colours)
#   ppt/theme/theme2.xml  -> used by the notes master  ("Office Theme" colours)
# The authored change swaps the two themes' colour schemes (and names) so that
# the slide master now carries the stock "Office Theme" palette and the notes
# master carries the old "Integral" palette.
#
# The PowerPoint object model only exposes the *presentation* theme (theme1.xml)
# for editing (via ThemeColorScheme / Master.ColorScheme) - there is no COM
# surface for the notes master's independent theme part - so we repaint the
# reachable theme (theme1.xml) with the "Office Theme" colour values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index -> (scheme slot, target "Office Theme" RGB as 0xBBGGRR for COM RGB())
$tcs.Colors(1).RGB  = 0x000000   # dk1      000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink 954F72
